$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date-formatted style from A49 onto A50:A59 (reuses existing
# style index / number format instead of Excel fabricating a brand new one).
$ws.Range("A49").Copy()
$ws.Range("A50:A59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Dates (all 2/9/2017 => serial 42775)
$ws.Range("A50").Value = 42775
$ws.Range("A51").Value = 42775
$ws.Range("A52").Value = 42775
$ws.Range("A53").Value = 42775
$ws.Range("A54").Value = 42775
$ws.Range("A55").Value = 42775
$ws.Range("A56").Value = 42775
$ws.Range("A57").Value = 42775
$ws.Range("A58").Value = 42775
$ws.Range("A59").Value = 42775

# "Who will it affect" column (E) -- enter first so "None" lands at the
# lowest new shared-string index, matching the source order.
$ws.Range("E50").Value = "None"
$ws.Range("E51").Value = "None"
$ws.Range("E52").Value = "None"
$ws.Range("E53").Value = "None"
$ws.Range("E54").Value = "None"
$ws.Range("E55").Value = "None"
$ws.Range("E56").Value = "None"
$ws.Range("E57").Value = "None"
$ws.Range("E58").Value = "None"
$ws.Range("E59").Value = "None"

# "Change / Action" column (B) -- entered out of row order in the source
# workbook (row 51 typed before row 50), reproduced here so the shared
# string table indices line up exactly.
$ws.Range("B51").Value = "Create Sales Laptops 1,2,3"
$ws.Range("B50").Value = "Create Sales Workstations 1,2,3,4,5"
$ws.Range("B52").Value = "Create Marketing Workstations 1,2,3"
$ws.Range("B53").Value = "Create Marketing Laptops 1,2,3"
$ws.Range("B54").Value = "Create Accounting Workstations 1,2,3,4,5"
$ws.Range("B55").Value = "Create Accounting Laptops 1,2"
$ws.Range("B56").Value = "Create Admins Workstations 1,2,3"
$ws.Range("B57").Value = "Create Admins Laptops 1,2,3"
$ws.Range("B58").Value = "Create Research Workstations 1,2,3,4"
$ws.Range("B59").Value = "Create Support Workstations"

# Remaining columns (C: Down Time?, D: How Long, F: Tech making changes,
# G: Approval) -- same recurring values used throughout the log.
foreach ($r in 50..59) {
    $ws.Range("C$r").Value = "No"
    $ws.Range("D$r").Value = "N/A"
    $ws.Range("F$r").Value = "Evan"
    $ws.Range("G$r").Value = "ES"
}

# Update the saved view state (scrolled / selected cell moved down).
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("A60").Select() | Out-Null
